$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# A new handoff package for "b.md" was generated (commit: "Generate Report
# for Handoff"). Update the Overview sheet and the two per-locale sheets
# (zh-cn, de-de) so row 3 (b.md) reflects the freshly generated handoff.
# ---------------------------------------------------------------------------

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1c6be98a4a47dfc1dd8ca875ff341ccd0309a5da/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e33ecc4b315b4a5a8489d639d4f8a4bbfbd882ef/e2e/b.md."

function Set-TextValue {
    param($range, [string]$text)
    # Plain ".Value = <word>" assignment auto-coerces text that looks like a
    # boolean (True/False, case-insensitive) into a real Boolean cell, which
    # is not what we want here (the column is plain text). Routing the
    # write through a formula + "paste values" keeps it a text cell.
    $range.Formula = '="' + $text + '"'
    $range.Copy() | Out-Null
    $range.PasteSpecial(-4163) | Out-Null   # xlPasteValues
}

# --- Overview sheet --------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-09-07 02:47:19"

# --- zh-cn sheet -------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
Set-TextValue $zhcn.Range("F3") "False"
$zhcn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("H3").Value = "2016-09-07 02:47:14"
$zhcn.Range("P3").Value = $errorDetail
$zhcn.Columns.Item(16).ColumnWidth = 39.15

# --- de-de sheet -------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
Set-TextValue $dede.Range("F3") "False"
$dede.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("H3").Value = "2016-09-07 02:47:19"
$dede.Range("P3").Value = $errorDetail
$dede.Columns.Item(16).ColumnWidth = 39.15

$excel.CutCopyMode = $false
